$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "C2" = 0.04898030544646303
    "D2" = 0.1257066524302282
    "E2" = 0.1035573562867231
    "F2" = 1.67036430150111
    "G2" = 0.002438876330124255
    "J2" = 0.1117533108220439
    "M2" = 1.352584820595027
    "N2" = 1.740294397118504
    "O2" = 4.574873944090655
    "C3" = 0.04348451995367952
    "D3" = 0.1260504724942564
    "E3" = 0.104899287764062
    "F3" = 1.640799753086426
    "G3" = 0.002443907155615908
    "J3" = 0.1140284108014908
    "M3" = 1.22724505357219
    "N3" = 1.606134505008129
    "O3" = 4.454729144343901
    "C4" = 0.04012880259830354
    "D4" = 0.126280195050569
    "E4" = 0.1057695588103202
    "F4" = 1.623870794883487
    "G4" = 0.002447160242067154
    "J4" = 0.1155032503235205
    "M4" = 1.150212736367109
    "N4" = 1.523993944342948
    "O4" = 4.384359006458055
    "C5" = 0.03876591023536946
    "D5" = 0.1263785061739817
    "E5" = 0.1061358538425878
    "F5" = 1.617278255967562
    "G5" = 0.002448527316004397
    "J5" = 0.1161238154003144
    "M5" = 1.118805339419026
    "N5" = 1.490583749881239
    "O5" = 4.356532256021751
    "C6" = 0.03853987680648174
    "D6" = 0.1263951148443585
    "E6" = 0.1061973807860754
    "F6" = 1.616202014522116
    "G6" = 0.002448756823119679
    "J6" = 0.1162280398046356
    "M6" = 1.113589260736376
    "N6" = 1.48503991737897
    "O6" = 4.3519628135071
    "C7" = 0.04011040372456876
    "D7" = 0.1262815018678385
    "E7" = 0.1057744516065631
    "F7" = 1.623780648117688
    "G7" = 0.002447178511114185
    "J7" = 0.1155115403504565
    "M7" = 1.149789227363186
    "N7" = 1.523543103533257
    "O7" = 4.383980291753971
    "C8" = 0.04708140169394426
    "D8" = 0.1258213476704988
    "E8" = 0.1040104437695435
    "F8" = 1.659915475881149
    "G8" = 0.002440576976344893
    "J8" = 0.1125215709533585
    "M8" = 1.309384196411017
    "N8" = 1.693989729170397
    "O8" = 4.532739360128971
    "C9" = 0.06090598405492642
    "D9" = 0.1250659768105287
    "E9" = 0.1009184843684068
    "F9" = 1.740563859347333
    "G9" = 0.002428927339792227
    "J9" = 0.1072780889896539
    "M9" = 1.621683565452898
    "N9" = 2.029942079240129
    "O9" = 4.851680228292707
    "C10" = 0.07116633368296732
    "D10" = 0.1245996481797214
    "E10" = 0.09887027245566848
    "F10" = 1.805900313395142
    "G10" = 0.002421149409398863
    "J10" = 0.1038056298581171
    "M10" = 1.850633252564307
    "N10" = 2.277632257103733
    "O10" = 5.102985925478947
    "C11" = 0.07585855701360344
    "D11" = 0.1244065461052912
    "E11" = 0.09798691408584226
    "F11" = 1.8369707072344
    "G11" = 0.002417778707945402
    "J11" = 0.1023088602583719
    "M11" = 1.954663130647134
    "N11" = 2.390466510659451
    "O11" = 5.221083733741636
    "C12" = 0.07763908447286383
    "D12" = 0.1243361436883177
    "E12" = 0.09765936181021395
    "F12" = 1.848932139035384
    "G12" = 0.002416526253252583
    "J12" = 0.101754030951863
    "M12" = 1.994037282722928
    "N12" = 2.433213383629663
    "O12" = 5.266353898608315
    "C13" = 0.07725545013198598
    "D13" = 0.1243511853361596
    "E13" = 0.09772959668354275
    "F13" = 1.846347295302508
    "G13" = 0.002416794928200891
    "J13" = 0.1018729904543338
    "M13" = 1.985558268557952
    "N13" = 2.424006295602794
    "O13" = 5.256579643447708
    "C14" = 0.07600496750835362
    "D14" = 0.1244006996107743
    "E14" = 0.09795982675307879
    "F14" = 1.837950847767758
    "G14" = 0.002417675188362432
    "J14" = 0.1022629740450931
    "M14" = 1.957902874285395
    "N14" = 2.393982962815755
    "O14" = 5.22479709900341
    "C15" = 0.07523949438964905
    "D15" = 0.1244313824400436
    "E15" = 0.09810175521136166
    "F15" = 1.832833328159808
    "G15" = 0.00241821748908606
    "J15" = 0.1025034098757249
    "M15" = 1.940960510433939
    "N15" = 2.3755951647193
    "O15" = 5.205401055879406
    "C16" = 0.07086019448175307
    "D16" = 0.1246126494989639
    "E16" = 0.09892897550383051
    "F16" = 1.803897081847452
    "G16" = 0.002421373052424875
    "J16" = 0.1039051189366322
    "M16" = 1.843832026891022
    "N16" = 2.270261138298622
    "O16" = 5.095344471309886
    "C17" = 0.06818006193844894
    "D17" = 0.124728714566956
    "E17" = 0.09944883969995377
    "F17" = 1.786492316169699
    "G17" = 0.002423351699304724
    "J17" = 0.1047862811262554
    "M17" = 1.784214287648325
    "N17" = 2.205680035018815
    "O17" = 5.028800010096859
    "C18" = 0.06664084145775462
    "D18" = 0.1247972645576922
    "E18" = 0.09975240598112478
    "F18" = 1.776608372821485
    "G18" = 0.002424505540240981
    "J18" = 0.1053008994612785
    "M18" = 1.749912540938794
    "N18" = 2.168549860450128
    "O18" = 4.990880604807955
    "C19" = 0.06612008261664926
    "D19" = 0.1248207827417893
    "E19" = 0.09985597052216555
    "F19" = 1.773283564942318
    "G19" = 0.002424898924448697
    "J19" = 0.1054764780897166
    "M19" = 1.738296708259369
    "N19" = 2.155980944334431
    "O19" = 4.978102575252478
    "C20" = 0.06846512550636419
    "D20" = 0.1247161738465792
    "E20" = 0.09939302794060101
    "F20" = 1.78833194484767
    "G20" = 0.002423139436979637
    "J20" = 0.1046916725835398
    "M20" = 1.790561874240353
    "N20" = 2.212553261873268
    "O20" = 5.035846976508083
    "C21" = 0.07637216355831811
    "D21" = 0.1243860823426814
    "E21" = 0.09789201381539936
    "F21" = 1.840411760676432
    "G21" = 0.002417415985497313
    "J21" = 0.1021481013010632
    "M21" = 1.96602648757505
    "N21" = 2.402801058347393
    "O21" = 5.234117453455042
    "C22" = 0.08156140485807839
    "D22" = 0.1241862020828179
    "E22" = 0.0969515609084084
    "F22" = 1.875590766485004
    "G22" = 0.002413814956259795
    "J22" = 0.1005554935542126
    "M22" = 2.080586815258073
    "N22" = 2.527247939251311
    "O22" = 5.366901954302136
    "C23" = 0.07878979641044737
    "D23" = 0.1242914365427197
    "E23" = 0.09744978864766285
    "F23" = 1.856709979456014
    "G23" = 0.00241572416557385
    "J23" = 0.1013990988491464
    "M23" = 2.019455199949903
    "N23" = 2.46081958113939
    "O23" = 5.295737306574949
    "C24" = 0.06833624322059961
    "D24" = 0.1247218378281261
    "E24" = 0.09941824583028502
    "F24" = 1.787499867959923
    "G24" = 0.002423235350120354
    "J24" = 0.1047344201089651
    "M24" = 1.787692213369894
    "N24" = 2.209445880822386
    "O24" = 5.032659990979027
    "C25" = 0.05714857791345196
    "D25" = 0.1252546883386039
    "E25" = 0.1017156698588266
    "F25" = 1.717686432853256
    "G25" = 0.002431941066315912
    "J25" = 0.1086300210937576
    "M25" = 1.985558268557952
    "N25" = 1.940960510433939
    "O25" = 5.256579643447708
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
